$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B14").Value = 3992
$ws.Range("G14").Value = 13857
$ws.Range("H14").Value = 5498
$ws.Range("L14").Value = 869
$ws.Range("N14").Value = 3068
$ws.Range("Q14").Value = 741
$ws.Range("S14").Value = 2226
$ws.Range("T14").Value = 7916
$ws.Range("U14").Value = 10772
$ws.Range("X14").Value = 5732
$ws.Range("AG14").Value = 8277
$ws.Range("AI14").Value = 113737
$ws.Range("B15").Value = 3821
$ws.Range("G15").Value = 13016
$ws.Range("H15").Value = 5298
$ws.Range("L15").Value = 665
$ws.Range("R15").Value = 2181
$ws.Range("T15").Value = 7478
$ws.Range("U15").Value = 10065
$ws.Range("X15").Value = 5081
$ws.Range("AI15").Value = 111987
$ws.Range("G16").Value = 13452
$ws.Range("H16").Value = 5214
$ws.Range("L16").Value = 696
$ws.Range("T16").Value = 7569
$ws.Range("U16").Value = 11532
$ws.Range("X16").Value = 5495
$ws.Range("AG16").Value = 9010
$ws.Range("AI16").Value = 118555
$ws.Range("G17").Value = 14489
$ws.Range("H17").Value = 5715
$ws.Range("R17").Value = 2541
$ws.Range("U17").Value = 12983
$ws.Range("Y17").Value = 4001
$ws.Range("AG17").Value = 9804
$ws.Range("AI17").Value = 125692
$ws.Range("B18").Value = 3940
$ws.Range("G18").Value = 15062
$ws.Range("H18").Value = 6050
$ws.Range("R18").Value = 2700
$ws.Range("T18").Value = 8551
$ws.Range("U18").Value = 13992
$ws.Range("X18").Value = 6314
$ws.Range("AG18").Value = 10536
$ws.Range("AI18").Value = 132514
$ws.Range("B19").Value = 4034
$ws.Range("G19").Value = 15355
$ws.Range("H19").Value = 6047
$ws.Range("I19").Value = 4007
$ws.Range("K19").Value = 377
$ws.Range("L19").Value = 793
$ws.Range("N19").Value = 3359
$ws.Range("Q19").Value = 899
$ws.Range("R19").Value = 2664
$ws.Range("T19").Value = 8989
$ws.Range("U19").Value = 15092
$ws.Range("V19").Value = 12540
$ws.Range("X19").Value = 6499
$ws.Range("Y19").Value = 4308
$ws.Range("AF19").Value = 126240
$ws.Range("AG19").Value = 11037
$ws.Range("AI19").Value = 137929
$ws.Range("G20").Value = 15215
$ws.Range("H20").Value = 5966
$ws.Range("I20").Value = 4035
$ws.Range("P20").Value = 2328
$ws.Range("R20").Value = 2480
$ws.Range("U20").Value = 15480
$ws.Range("V20").Value = 12832
$ws.Range("AF20").Value = 128464
$ws.Range("AG20").Value = 11227
$ws.Range("AI20").Value = 140306
$ws.Range("G21").Value = 15221
$ws.Range("H21").Value = 6068
$ws.Range("I21").Value = 3921
$ws.Range("Q21").Value = 842
$ws.Range("R21").Value = 2494
$ws.Range("T21").Value = 9192
$ws.Range("U21").Value = 15803
$ws.Range("V21").Value = 13055
$ws.Range("AF21").Value = 131335
$ws.Range("AG21").Value = 11505
$ws.Range("AH21").Value = 634
$ws.Range("AI21").Value = 143475
$ws.Range("G22").Value = 15036
$ws.Range("H22").Value = 5877
$ws.Range("I22").Value = 3911
$ws.Range("N22").Value = 3401
$ws.Range("Q22").Value = 898
$ws.Range("R22").Value = 2464
$ws.Range("T22").Value = 9512
$ws.Range("U22").Value = 16025
$ws.Range("V22").Value = 13204
$ws.Range("X22").Value = 7315
$ws.Range("AF22").Value = 133409
$ws.Range("AG22").Value = 11747
$ws.Range("AI22").Value = 145801
$ws.Range("B23").Value = 4368
$ws.Range("G23").Value = 15302
$ws.Range("H23").Value = 6044
$ws.Range("I23").Value = 4163
$ws.Range("L23").Value = 839
$ws.Range("N23").Value = 3441
$ws.Range("P23").Value = 2195
$ws.Range("Q23").Value = 842
$ws.Range("R23").Value = 2576
$ws.Range("U23").Value = 16606
$ws.Range("V23").Value = 13722
$ws.Range("X23").Value = 7407
$ws.Range("Y23").Value = 5015
$ws.Range("AF23").Value = 134979
$ws.Range("AG23").Value = 12176
$ws.Range("AH23").Value = 696
$ws.Range("AI23").Value = 147852
$ws.Range("B24").Value = 4451
$ws.Range("E24").Value = 14108
$ws.Range("G24").Value = 15892
$ws.Range("H24").Value = 6371
$ws.Range("I24").Value = 4243
$ws.Range("J24").Value = 2141
$ws.Range("N24").Value = 3551
$ws.Range("Q24").Value = 874
$ws.Range("R24").Value = 2617
$ws.Range("S24").Value = 4255
$ws.Range("U24").Value = 17203
$ws.Range("V24").Value = 14237
$ws.Range("X24").Value = 7586
$ws.Range("Y24").Value = 5374
$ws.Range("AF24").Value = 139999
$ws.Range("AG24").Value = 12672
$ws.Range("AH24").Value = 748
$ws.Range("AI24").Value = 153419
$ws.Range("B25").Value = 4387
$ws.Range("G25").Value = 15921
$ws.Range("H25").Value = 6292
$ws.Range("I25").Value = 4239
$ws.Range("J25").Value = 2062
$ws.Range("K25").Value = 367
$ws.Range("N25").Value = 3690
$ws.Range("R25").Value = 2705
$ws.Range("T25").Value = 9725
$ws.Range("U25").Value = 17234
$ws.Range("V25").Value = 14230
$ws.Range("X25").Value = 7742
$ws.Range("AF25").Value = 141285
$ws.Range("AG25").Value = 12809
$ws.Range("AI25").Value = 154766
$ws.Range("B26").Value = 4318
$ws.Range("G26").Value = 15334
$ws.Range("H26").Value = 6163
$ws.Range("I26").Value = 4323
$ws.Range("K26").Value = 290
$ws.Range("L26").Value = 801
$ws.Range("N26").Value = 3427
$ws.Range("P26").Value = 2413
$ws.Range("Q26").Value = 827
$ws.Range("R26").Value = 2640
$ws.Range("S26").Value = 4153
$ws.Range("T26").Value = 8338
$ws.Range("U26").Value = 15797
$ws.Range("V26").Value = 13861
$ws.Range("X26").Value = 6362
$ws.Range("AF26").Value = 132889
$ws.Range("AG26").Value = 11972
$ws.Range("AI26").Value = 145498
